$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we are about to rewrite keep their
# text formatting so numeric-looking strings (e.g. "1.000", "75.40")
# are not coerced into numbers by the Value setter.
$ws.Range("D2:D15").NumberFormat = "@"
$ws.Range("D17:D19").NumberFormat = "@"
$ws.Range("D21:D30").NumberFormat = "@"
$ws.Range("D32:D50").NumberFormat = "@"

$ws.Range("D2").Value = "25.581.89"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "1.669.61"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "237.97"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4773"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D8").Value = "0.2616"
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("D9").Value = "0.06176"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").Value = "1.671.12"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").Value = "0.06990"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "14.83"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "0.5891"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("D14").Value = "4.376"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "75.40"
$ws.Range("E15").Value = "  +4.19%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "0.9996"
$ws.Range("D18").Value = "25.572.77"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").Value = "0.000006747"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("D21").Value = "1.885.31"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("D22").Value = "4.444"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").Value = "8.801"
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("D24").Value = "5.265"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "136.78"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").Value = "15.05"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("D27").Value = "1.384"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "1.720"
$ws.Range("E28").Value = "  +5.56%  "
$ws.Range("D29").Value = "104.73"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "3.995"
$ws.Range("E30").Value = "  +6.78%  "
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").Value = "3.635"
$ws.Range("D33").Value = "0.9990"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "0.04306"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").Value = "2.621"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "0.9556"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("D37").Value = "0.6049"
$ws.Range("E37").Value = "  +5.32%  "
$ws.Range("D38").Value = "2.576"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").Value = "0.9231"
$ws.Range("E39").Value = "  +12.58%  "
$ws.Range("D40").Value = "0.9997"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "1.853"
$ws.Range("E41").Value = "  +4.50%  "
$ws.Range("D42").Value = "0.01473"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").Value = "97.76"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "0.3761"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "4.876"
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").Value = "0.1119"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").Value = "6.218"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("D48").Value = "0.05269"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "29.96"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "7.442"
$ws.Range("E50").Value = "  +4.34%  "
$ws.Range("E51").Value = "  +0.18%  "
